$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 52
$newRow = 53

# Carry the formatting (incl. the A-column date style) down from the
# preceding row, the same way dragging the fill handle would in Excel,
# so the new row reuses the existing style instead of minting a new one.
$ws.Range("A" + $lastRow + ":F" + $lastRow).Copy() | Out-Null
$ws.Range("A" + $newRow + ":F" + $newRow).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($newRow, 1).Value = 45383
$ws.Cells.Item($newRow, 2).Value = -0.406
$ws.Cells.Item($newRow, 3).Value = 0.562
$ws.Cells.Item($newRow, 4).Value = -0.649
$ws.Cells.Item($newRow, 5).Value = 0.436
$ws.Cells.Item($newRow, 6).Value = 1.657
